# Changes for RAD, VLink and Multibill
#
# Katalon test-run timestamps in the "Date" column (column B) of the
# PersonalSearchTransactions and ConfirmAcNumError sheets are refreshed
# to reflect the latest test execution run.

$wb = $excel.ActiveWorkbook

$personalTimestamps = @(
    "Fri Nov 14 21:58:42 EST 2025",
    "Fri Nov 14 21:59:25 EST 2025",
    "Fri Nov 14 22:00:08 EST 2025",
    "Fri Nov 14 22:00:55 EST 2025",
    "Fri Nov 14 22:01:42 EST 2025",
    "Fri Nov 14 22:02:23 EST 2025",
    "Fri Nov 14 22:03:07 EST 2025",
    "Fri Nov 14 22:03:49 EST 2025",
    "Fri Nov 14 22:04:33 EST 2025",
    "Fri Nov 14 22:05:17 EST 2025",
    "Fri Nov 14 22:06:00 EST 2025",
    "Fri Nov 14 22:06:42 EST 2025",
    "Fri Nov 14 22:07:25 EST 2025",
    "Fri Nov 14 22:08:08 EST 2025",
    "Fri Nov 14 22:08:50 EST 2025",
    "Fri Nov 14 22:09:33 EST 2025",
    "Fri Nov 14 22:10:15 EST 2025"
)

$confirmAcNumTimestamps = @(
    "Fri Nov 14 22:10:57 EST 2025",
    "Fri Nov 14 22:11:52 EST 2025",
    "Fri Nov 14 22:12:46 EST 2025",
    "Fri Nov 14 22:13:44 EST 2025",
    "Fri Nov 14 22:14:37 EST 2025",
    "Fri Nov 14 22:15:32 EST 2025",
    "Fri Nov 14 22:16:26 EST 2025",
    "Fri Nov 14 22:17:21 EST 2025",
    "Fri Nov 14 22:18:15 EST 2025",
    "Fri Nov 14 22:19:10 EST 2025",
    "Fri Nov 14 22:20:06 EST 2025",
    "Fri Nov 14 22:21:00 EST 2025",
    "Fri Nov 14 22:21:54 EST 2025",
    "Fri Nov 14 22:22:50 EST 2025",
    "Fri Nov 14 22:23:44 EST 2025",
    "Fri Nov 14 22:24:39 EST 2025",
    "Fri Nov 14 22:25:33 EST 2025"
)

$saleSearchSheet = $wb.Worksheets.Item("SaleSearchTransactions")
for ($i = 0; $i -lt $personalTimestamps.Length; $i++) {
    $row = $i + 2
    $saleSearchSheet.Cells.Item($row, 2).Value = $personalTimestamps[$i]
}

$personalSheet = $wb.Worksheets.Item("PersonalSearchTransactions")
for ($i = 0; $i -lt $confirmAcNumTimestamps.Length; $i++) {
    $row = $i + 2
    $personalSheet.Cells.Item($row, 2).Value = $confirmAcNumTimestamps[$i]
}
